# Update countries & provincias Spain
# - Refresh case counts for Estados Unidos (row 4), Brasil (row 14), Nigeria (row 85)
# - Refresh Mayotte stats and swap its row position with Mali (rows 116/117)
# - Bump the "Datos actualizados" timestamp from 00:52 to 01:22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 01:22"

# Estados Unidos - row 4
$ws.Range("B4").Value = 957356
$ws.Range("C4").Value = 32124
$ws.Range("D4").Value = 116201
$ws.Range("E4").Value = 787011
$ws.Range("G4").Value = 1951
$ws.Range("H4").Value = 54144

# Brasil - row 14
$ws.Range("B14").Value = 58997
$ws.Range("C14").Value = 6002
$ws.Range("E14").Value = 25800
$ws.Range("G14").Value = 367
$ws.Range("H14").Value = 4037

# Nigeria - row 85
$ws.Range("B85").Value = 1182
$ws.Range("C85").Value = 87
$ws.Range("D85").Value = 222
$ws.Range("E85").Value = 925
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 35

# Mayotte now occupies row 116 (previously Mali), with refreshed totals
$ws.Range("A116").Value = "Mayotte"
$ws.Range("B116").Value = 380
$ws.Range("C116").Value = 26
$ws.Range("D116").Value = 144
$ws.Range("E116").Value = 232
$ws.Range("F116").Value = 4
$ws.Range("H116").Value = 4

# Mali now occupies row 117 (previously Mayotte), keeping its prior totals
$ws.Range("A117").Value = "Mali"
$ws.Range("B117").Value = 370
$ws.Range("C117").Value = 45
$ws.Range("D117").Value = 91
$ws.Range("E117").Value = 258
$ws.Range("F117").Value = 0
$ws.Range("H117").Value = 21
